$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2:T11").ClearContents()

# Column A
$ws.Range("A2").Value = 'M1'
$ws.Range("A3").Value = 'M1'
$ws.Range("A4").Value = 'M1'
$ws.Range("A5").Value = 'M1'
$ws.Range("A6").Value = 'M1'
$ws.Range("A7").Value = 'M2'
$ws.Range("A8").Value = 'M2'
$ws.Range("A9").Value = 'M2'
$ws.Range("A10").Value = 'M2'
$ws.Range("A11").Value = 'M2'

# Column B
$ws.Range("B2").Value = 'Matn1'
$ws.Range("B3").Value = 'Matn1'
$ws.Range("B4").Value = 'Matn1'
$ws.Range("B5").Value = 'Matn1'
$ws.Range("B6").Value = 'Matn1'
$ws.Range("B7").Value = 'Matn1'
$ws.Range("B8").Value = 'Matn1'
$ws.Range("B9").Value = 'Matn1'
$ws.Range("B10").Value = 'Matn1'
$ws.Range("B11").Value = 'Matn1'

# Column C
$ws.Range("C2").Value = 'Itga1'
$ws.Range("C3").Value = 'Itga1'
$ws.Range("C4").Value = 'Itga1'
$ws.Range("C5").Value = 'Itga1'
$ws.Range("C6").Value = 'Itga1'
$ws.Range("C7").Value = 'Itga1'
$ws.Range("C8").Value = 'Itga1'
$ws.Range("C9").Value = 'Itga1'
$ws.Range("C10").Value = 'Itga1'
$ws.Range("C11").Value = 'Itga1'

# Column D
$ws.Range("D2").Value = 'ECs'
$ws.Range("D3").Value = 'FAPs'
$ws.Range("D4").Value = 'M1'
$ws.Range("D5").Value = 'M2'
$ws.Range("D6").Value = 'sCs'
$ws.Range("D7").Value = 'ECs'
$ws.Range("D8").Value = 'FAPs'
$ws.Range("D9").Value = 'M1'
$ws.Range("D10").Value = 'M2'
$ws.Range("D11").Value = 'sCs'

# Numeric columns
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2053186666666666
$ws.Range("H2").Value = 0.6159559999999999
$ws.Range("I2").Value = 0.420847658662702
$ws.Range("J2").Value = 0.4208476586627021
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 51.78202533333334
$ws.Range("N2").Value = 155.346076
$ws.Range("O2").Value = 0.7202935140152373
$ws.Range("P2").Value = 0.7202935140152373
$ws.Range("Q2").Value = 10.63181639873956
$ws.Range("R2").Value = 95.68634758865601
$ws.Range("S2").Value = 0.3031338389232427
$ws.Range("T2").Value = 0.3031338389232428
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2053186666666666
$ws.Range("H3").Value = 0.6159559999999999
$ws.Range("I3").Value = 0.420847658662702
$ws.Range("J3").Value = 0.4208476586627021
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.38032666666667
$ws.Range("N3").Value = 31.14098
$ws.Range("O3").Value = 0.1443914548190984
$ws.Range("P3").Value = 0.1443914548190984
$ws.Range("Q3").Value = 2.131274830764444
$ws.Range("R3").Value = 19.18147347688
$ws.Range("S3").Value = 0.06076680569151887
$ws.Range("T3").Value = 0.06076680569151888
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2053186666666666
$ws.Range("H4").Value = 0.6159559999999999
$ws.Range("I4").Value = 0.420847658662702
$ws.Range("J4").Value = 0.4208476586627021
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.523303
$ws.Range("N4").Value = 1.569909
$ws.Range("O4").Value = 0.007279200733040383
$ws.Range("P4").Value = 0.007279200733040383
$ws.Range("Q4").Value = 0.1074438742226666
$ws.Range("R4").Value = 0.9669948680039999
$ws.Range("S4").Value = 0.003063434585435869
$ws.Range("T4").Value = 0.00306343458543587
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.2053186666666666
$ws.Range("H5").Value = 0.6159559999999999
$ws.Range("I5").Value = 0.420847658662702
$ws.Range("J5").Value = 0.4208476586627021
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.511976
$ws.Range("N5").Value = 10.535928
$ws.Range("O5").Value = 0.04885196200598933
$ws.Range("P5").Value = 0.04885196200598933
$ws.Range("Q5").Value = 0.7210742296853333
$ws.Range("R5").Value = 6.489668067167999
$ws.Range("S5").Value = 0.02055923383129989
$ws.Range("T5").Value = 0.02055923383129989
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.2053186666666666
$ws.Range("H6").Value = 0.6159559999999999
$ws.Range("I6").Value = 0.420847658662702
$ws.Range("J6").Value = 0.4208476586627021
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.692542
$ws.Range("N6").Value = 17.077626
$ws.Range("O6").Value = 0.07918386842663461
$ws.Range("P6").Value = 0.07918386842663462
$ws.Range("Q6").Value = 1.168785133384
$ws.Range("R6").Value = 10.519066200456
$ws.Range("S6").Value = 0.03332434563120463
$ws.Range("T6").Value = 0.03332434563120464
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.2825506666666667
$ws.Range("H7").Value = 0.847652
$ws.Range("I7").Value = 0.579152341337298
$ws.Range("J7").Value = 0.579152341337298
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 51.78202533333334
$ws.Range("N7").Value = 155.346076
$ws.Range("O7").Value = 0.7202935140152373
$ws.Range("P7").Value = 0.7202935140152373
$ws.Range("Q7").Value = 14.63104577928356
$ws.Range("R7").Value = 131.679412013552
$ws.Range("S7").Value = 0.4171596750919946
$ws.Range("T7").Value = 0.4171596750919946
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.2825506666666667
$ws.Range("H8").Value = 0.847652
$ws.Range("I8").Value = 0.579152341337298
$ws.Range("J8").Value = 0.579152341337298
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 10.38032666666667
$ws.Range("N8").Value = 31.14098
$ws.Range("O8").Value = 0.1443914548190984
$ws.Range("P8").Value = 0.1443914548190984
$ws.Range("Q8").Value = 2.932968219884445
$ws.Range("R8").Value = 26.39671397896
$ws.Range("S8").Value = 0.0836246491275795
$ws.Range("T8").Value = 0.0836246491275795
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.2825506666666667
$ws.Range("H9").Value = 0.847652
$ws.Range("I9").Value = 0.579152341337298
$ws.Range("J9").Value = 0.579152341337298
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.523303
$ws.Range("N9").Value = 1.569909
$ws.Range("O9").Value = 0.007279200733040383
$ws.Range("P9").Value = 0.007279200733040383
$ws.Range("Q9").Value = 0.1478596115186667
$ws.Range("R9").Value = 1.330736503668
$ws.Range("S9").Value = 0.004215766147604514
$ws.Range("T9").Value = 0.004215766147604514
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2825506666666667
$ws.Range("H10").Value = 0.847652
$ws.Range("I10").Value = 0.579152341337298
$ws.Range("J10").Value = 0.579152341337298
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.511976
$ws.Range("N10").Value = 10.535928
$ws.Range("O10").Value = 0.04885196200598933
$ws.Range("P10").Value = 0.04885196200598933
$ws.Range("Q10").Value = 0.9923111601173334
$ws.Range("R10").Value = 8.930800441056
$ws.Range("S10").Value = 0.02829272817468945
$ws.Range("T10").Value = 0.02829272817468944
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.2825506666666667
$ws.Range("H11").Value = 0.847652
$ws.Range("I11").Value = 0.579152341337298
$ws.Range("J11").Value = 0.579152341337298
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 5.692542
$ws.Range("N11").Value = 17.077626
$ws.Range("O11").Value = 0.07918386842663461
$ws.Range("P11").Value = 0.07918386842663462
$ws.Range("Q11").Value = 1.608431537128
$ws.Range("R11").Value = 14.475883834152
$ws.Range("S11").Value = 0.04585952279542998
$ws.Range("T11").Value = 0.04585952279542999
